$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Archipel B.V."
$ws.Cells.Item($row, 3).Value = "KDV"

# "2024-04-02" looks like a date, so Excel would normally auto-convert it to a
# date serial number. Force it to stay plain text (matching the source data,
# which stores every date in this column as a literal string), then restore
# the cell's style so no extra formatting is left behind.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-04-02"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
